$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.127.39'
$ws.Range("E2").Value = '  -4.40%  '

$ws.Range("D3").Value = '1.652.01'
$ws.Range("E3").Value = '  -3.44%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.89'
$ws.Range("E5").Value = '  -3.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5117'
$ws.Range("E6").Value = '  -3.35%  '

$ws.Range("E7").Value = '  +0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2585'
$ws.Range("E8").Value = '  -2.95%  '

$ws.Range("E9").Value = '  -3.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.00'
$ws.Range("E10").Value = '  -4.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07791'
$ws.Range("E11").Value = '  +1.26%  '

$ws.Range("D12").Value = '1.652.86'
$ws.Range("E12").Value = '  -3.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.285'
$ws.Range("E13").Value = '  -4.70%  '

$ws.Range("D14").Value = '1.880.61'
$ws.Range("E14").Value = '  -3.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5523'
$ws.Range("E15").Value = '  -5.25%  '

$ws.Range("D16").Value = '0.0₅8035'
$ws.Range("E16").Value = '  -2.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.16'
$ws.Range("E17").Value = '  -5.68%  '

$ws.Range("D18").Value = '26.163.60'
$ws.Range("E18").Value = '  -4.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.33'
$ws.Range("E20").Value = '  -5.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.410'
$ws.Range("E21").Value = '  -4.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.07'
$ws.Range("E22").Value = '  -3.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.047'
$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("E24").Value = '  +0.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.81'
$ws.Range("E25").Value = '  -0.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.743'
$ws.Range("E26").Value = '  +3.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1183'
$ws.Range("E27").Value = '  -1.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.988'
$ws.Range("E28").Value = '  -3.32%  '

$ws.Range("E29").Value = '  -2.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05132'
$ws.Range("E30").Value = '  -3.71%  '

$ws.Range("E31").Value = '  -3.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.344'
$ws.Range("E32").Value = '  -3.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.221'
$ws.Range("E33").Value = '  -6.22%  '

$ws.Range("E34").Value = '  -4.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.748'
$ws.Range("E35").Value = '  -4.46%  '

$ws.Range("E36").Value = '  -2.61%  '

$ws.Range("E37").Value = '  -1.32%  '

$ws.Range("D38").Value = '1.165.30'
$ws.Range("E38").Value = '  +2.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5696'
$ws.Range("E39").Value = '  -2.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01586'
$ws.Range("E40").Value = '  -3.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.559'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("E42").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8300'
$ws.Range("E43").Value = '  -1.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.654'
$ws.Range("E44").Value = '  -2.53%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.12'
$ws.Range("E45").Value = '  -1.10%  '

$ws.Range("D46").Value = '1.790.66'
$ws.Range("E46").Value = '  -3.40%  '

$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4554'
$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.63'
$ws.Range("E49").Value = '  -3.54%  '

$ws.Range("E50").Value = '  +0.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.887'
$ws.Range("E51").Value = '  -2.72%  '
